# Update the ObjTables header metadata (version bump 0.0.9 -> 1.0.0, new timestamp)
# on every sheet that carries one, and backfill the "Verbose name" column of the
# Schema table for every attribute row.

$wb = $excel.ActiveWorkbook

foreach ($sheet in $wb.Worksheets) {
    $sheet.Unprotect()
}

$newVersion = "1.0.0"
$newDate = "2020-05-29 00:18:51"

$toc = $wb.Worksheets.Item("!!_Table of contents")
$toc.Range("A1").Value = "!!!ObjTables objTablesVersion='$newVersion' date='$newDate'"
$toc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='$newDate' objTablesVersion='$newVersion'"

$schema = $wb.Worksheets.Item("!!_Schema")
$schema.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='$newDate' objTablesVersion='$newVersion'"

$company = $wb.Worksheets.Item("!!Company")
$company.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' class='Company' name='Companies' date='$newDate' objTablesVersion='$newVersion'"

$people = $wb.Worksheets.Item("!!People")
$people.Range("A1").Value = "!!ObjTables type='Data' tableFormat='row' class='Person' name='People' date='$newDate' objTablesVersion='$newVersion'"

# Fill in the "Verbose name" (column E) for every attribute row on the Schema sheet.
$verboseNames = @{
    4  = "Address"
    5  = "Name"
    6  = "URL"
    8  = "Address"
    9  = "Company"
    10 = "Email address"
    11 = "Name"
    12 = "Phone number"
    13 = "Type"
    15 = "City"
    16 = "Country"
    17 = "State"
    18 = "Street"
    19 = "Zip code"
}

foreach ($row in $verboseNames.Keys) {
    $schema.Cells.Item($row, 5).Value = $verboseNames[$row]
}

foreach ($sheet in $wb.Worksheets) {
    $sheet.Protect()
}
